$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Paste in the newly-uploaded daily rows (45992 .. 46002 / rows 35-45):
#    date (col A) was already present - fill in total impressions (B),
#    ad impressions (C), visitors (E), inquiries (F), receptions (G),
#    ad spend (H), leads (I) and deal amount (J). Column D (delta) is a
#    formula, added separately below so the B17:D63 shared-formula chain
#    stays intact.
# ---------------------------------------------------------------------------
$rowData = @(
    @(35, 1740, 1549, 55, 6, 7, 227, 7, 128),
    @(36, 1866, 1609, 55, 5, 6, 287, 4, 36),
    @(37, 1392, 1609, 46, 5, 6, 206, 13, 32),
    @(38, 1364, 1084, 38, 3, 6, 209, 4, 61),
    @(39, 2956, 2689, 32, 3, 4, 278, 4, 0),
    @(40, 1590, 1409, 22, 3, 3, 0, 1, 0),
    @(41, 2060, 1939, 32, 1, 1, 245, 5, 96),
    @(42, 3378, 3186, 67, 8, 8, 314, 12, 0),
    @(43, 2722, 2438, 59, 6, 6, 290, 8, 83),
    @(44, 1974, 1703, 37, 5, 5, 222, 10, 115),
    @(45, 2660, 2176, 61, 11, 11, 278, 16, 36)
)

foreach ($r in $rowData) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value  = $r[1]   # B - total impressions
    $ws.Cells.Item($row, 3).Value  = $r[2]   # C - ad impressions
    $ws.Cells.Item($row, 5).Value  = $r[3]   # E - visitors
    $ws.Cells.Item($row, 6).Value  = $r[4]   # F - inquiries
    $ws.Cells.Item($row, 7).Value  = $r[5]   # G - receptions
    $ws.Cells.Item($row, 8).Value  = $r[6]   # H - ad spend
    $ws.Cells.Item($row, 9).Value  = $r[7]   # I - leads
    $ws.Cells.Item($row, 10).Value = $r[8]   # J - deal amount
}

# ---------------------------------------------------------------------------
# 2) Extend the "B-C" delta formula in column D all the way down to row 63
#    (it previously stopped at row 34 with a few rows left as plain, non
#    shared formulas). Re-enter every cell from D17 to D63 so the whole
#    block becomes one consistent formula run.
# ---------------------------------------------------------------------------
for ($row = 17; $row -le 63; $row++) {
    $ws.Cells.Item($row, 4).Formula = "=B$row-C$row"
}

# ---------------------------------------------------------------------------
# 3) View state: scroll the frozen-pane window so B3 is the top-left cell
#    of the scrollable area, and move the active selection to M8.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 2
$ws.Range("M8").Select() | Out-Null
